# Update cryptocurrency price/volume data (Price = column D, Volume(1h) = column E).
# Column D mixes "text that looks numeric" (e.g. "1.001") with multi-dot
# strings that can never be numbers (e.g. "28.598.46"). Excel auto-coerces
# numeric-looking text typed into a cell into a real number, which would
# flip the cell's stored type away from the original text/string type. To
# preserve the original text representation (leading/trailing zeros, dot
# groupings, etc.) we force the cell to Text format ("@") immediately
# before writing any value that Excel could otherwise reinterpret as a
# number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.598.46"
$ws.Range("E2").Value = "  +0.84%  "

$ws.Range("D3").Value = "1.802.62"
$ws.Range("E3").Value = "  -0.56%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.94"
$ws.Range("E5").Value = "  -0.23%  "

$ws.Range("E6").Value = "  +0.08%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5466"
$ws.Range("E7").Value = "  -4.62%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3768"
$ws.Range("E8").Value = "  -2.80%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07491"
$ws.Range("E9").Value = "  -1.60%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.36"
$ws.Range("E10").Value = "  -1.78%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.113"
$ws.Range("E11").Value = "  -2.29%  "

$ws.Range("E12").Value = "  +0.10%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.62"
$ws.Range("E13").Value = "  -2.89%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.148"
$ws.Range("E14").Value = "  -1.72%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.386"
$ws.Range("E15").Value = "  +1.61%  "

$ws.Range("D16").Value = "1.797.25"
$ws.Range("E16").Value = "  -0.70%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "90.20"
$ws.Range("E17").Value = "  -2.02%  "

$ws.Range("E18").Value = "  -1.14%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06453"
$ws.Range("E19").Value = "  -0.36%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.24"
$ws.Range("E21").Value = "  -0.54%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.919"

$ws.Range("D23").Value = "28.611.60"
$ws.Range("E23").Value = "  +0.84%  "

$ws.Range("E24").Value = "  -1.76%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.091"
$ws.Range("E25").Value = "  -1.71%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "158.74"
$ws.Range("E26").Value = "  +0.27%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.44"
$ws.Range("E27").Value = "  -1.94%  "

$ws.Range("D28").Value = "2.004.97"
$ws.Range("E28").Value = "  -0.77%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.350"
$ws.Range("E29").Value = "  -3.50%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "122.87"
$ws.Range("E30").Value = "  -0.81%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.105"
$ws.Range("E31").Value = "  -5.21%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1061"
$ws.Range("E32").Value = "  +0.29%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.640"
$ws.Range("E33").Value = "  -2.39%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.682"
$ws.Range("E34").Value = "  +1.26%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.06497"
$ws.Range("E35").Value = "  +6.56%  "

$ws.Range("E36").Value = "  +3.71%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02301"
$ws.Range("E37").Value = "  -0.93%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.730"
$ws.Range("E38").Value = "  -2.62%  "

$ws.Range("E39").Value = "  -0.35%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "11.24"
$ws.Range("E40").Value = "  -3.77%  "

$ws.Range("E41").Value = "  +3.45%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6227"
$ws.Range("E42").Value = "  -3.20%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.432"
$ws.Range("E43").Value = "  +4.04%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.000"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.28"
$ws.Range("E45").Value = "  -0.86%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.691"
$ws.Range("E46").Value = "  -0.37%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5841"
$ws.Range("E47").Value = "  -2.72%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "126.51"
$ws.Range("E48").Value = "  +3.12%  "

$ws.Range("E49").Value = "  -0.41%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.156"
$ws.Range("E50").Value = "  +0.58%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06889"
$ws.Range("E51").Value = "  +0.58%  "
